$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate text in A1 -------------------------
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.06 = 24169.7 pesos`n✅ 24169.7 pesos = 6.06 = 951.04 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10 / O10 / N12 / O12 values -----------------
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 165
$ws2.Range("O10").Value = 3988
$ws2.Range("N12").Value = 3990
$ws2.Range("O12").Value = 157
